$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.64
$ws.Range("W2").Value = 2.56

$ws.Range("G3").Value = 9
$ws.Range("J3").Value = 5.2
$ws.Range("N3").Value = 5
$ws.Range("W3").Value = 1.14

$ws.Range("F4").Value = 1.28
$ws.Range("G4").Value = 1.67
$ws.Range("H4").Value = 1.2
$ws.Range("J4").Value = 1.4
$ws.Range("M4").Value = 1.03
$ws.Range("S4").Value = 1.67
$ws.Range("T4").Value = 1.94
$ws.Range("W4").Value = 2.4

$ws.Range("F5").Value = 2.16
$ws.Range("G5").Value = 2.4
$ws.Range("H5").Value = 3.45
$ws.Range("I5").Value = 4.1
$ws.Range("J5").Value = 3.15
$ws.Range("K5").Value = 3.7
$ws.Range("M5").Value = 1.08
$ws.Range("N5").Value = 3.05
$ws.Range("O5").Value = 1.32
$ws.Range("P5").Value = 1.73
$ws.Range("Q5").Value = 2.1
$ws.Range("R5").Value = 1.24
$ws.Range("S5").Value = 3.45
$ws.Range("T5").Value = 1.85
$ws.Range("U5").Value = 1.97
$ws.Range("V5").Value = 1.32
$ws.Range("W5").Value = 1.71
$ws.Range("X5").Value = 14.5
$ws.Range("Y5").Value = 15
$ws.Range("Z5").Value = 32
$ws.Range("AB5").Value = 11
$ws.Range("AC5").Value = 9.199999999999999
$ws.Range("AD5").Value = 19
$ws.Range("AE5").Value = 60
$ws.Range("AF5").Value = 17
$ws.Range("AG5").Value = 13.5
$ws.Range("AH5").Value = 23
$ws.Range("AJ5").Value = 38
$ws.Range("AK5").Value = 32
$ws.Range("AL5").Value = 55
$ws.Range("AN5").Value = 27
$ws.Range("AO5").Value = 70

$ws.Range("F6").Value = 1.78
$ws.Range("G6").Value = 1.88
$ws.Range("K6").Value = 3.7
$ws.Range("L6").Value = 1.46
$ws.Range("W6").Value = 2.08
$ws.Range("Z6").Value = 1000
$ws.Range("AD6").Value = 990
$ws.Range("AJ6").Value = 1000
$ws.Range("AN6").Value = 20
